# Add incident support gateway
# - Rename "SystemID" variable key to lowercase "systemID" (A8)
# - Add a new row 17 for a new variable "incidentSupported"
# - Update the sheet's selection/scroll state to reflect the edit position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the variable name casing in row 8 (SystemID -> systemID)
$ws.Range("A8").Value = "systemID"

# Copy formatting (row height, styles, borders, wrap text, etc.) from row 16
# down to the new row 17 so the new row matches the look of the existing table
$ws.Range("A16:D16").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)

# Populate the new "incidentSupported" variable row
$ws.Range("A17").Value = "incidentSupported"
$ws.Range("B17").Value = "Is this incident supported?"
$ws.Range("C17").Value = "Boolean"
$ws.Range("D17").Value = "Company"

$excel.CutCopyMode = 0

# Update the view so the new row is visible/selected, similar to the author's
# final cursor position after adding the row
$ws.Activate()
$ws.Range("B18").Select()
